$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The attendance roll-over crossed midnight: 15/11/2024 22:43 -> 16/11/2024 00:24,
# so the day-16 column (S) now holds what was being recorded, and day-15 (R) is cleared.
$ws.Range("C6").Value = "16/11/2024 00:24"

# Row 8 (student 1): was marked Present ("P") on day 15 (R8); same mark now lands
# on day 16 (S8). Reuse R8's exact cell format (green "Present" fill) for S8.
$ws.Range("R8").Copy()
$ws.Range("S8").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("S8").Value = "P"

# Rows 9-32 (students 2-25): all recorded as Absent ("F", red fill) on day 16.
# Reuse R10's exact format (the existing red "Absent" style already in the sheet).
$ws.Range("R10").Copy()
$ws.Range("S9:S32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("S9:S32").Value = "F"

# Day-15 column (R) is now blank for all of these rows; restore the plain
# unmarked-cell format (matching the neighbouring blank attendance cells).
$ws.Range("Q8").Copy()
$ws.Range("R8:R32").PasteSpecial(-4122)
$excel.CutCopyMode = $false
$ws.Range("R8:R32").ClearContents()
